$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current extr1..extr8 rows (rows 8..15) before overwriting anything ---
$vals = @{}
for ($r = 8; $r -le 15; $r++) {
    $vals[$r] = @(
        $ws.Cells.Item($r, 1).Value(),
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value(),
        $ws.Cells.Item($r, 4).Value(),
        $ws.Cells.Item($r, 5).Value()
    )
}

# --- Step 2: write extr1..extr8 into their new home, rows 10..17 (shifted down by 2) ---
for ($r = 8; $r -le 15; $r++) {
    $destRow = $r + 2
    $row = $vals[$r]
    $ws.Cells.Item($destRow, 1).Value = $row[0]
    $ws.Cells.Item($destRow, 2).Value = $row[1]
    $ws.Cells.Item($destRow, 3).Value = $row[2]
    $ws.Cells.Item($destRow, 4).Value = $row[3]
    $ws.Cells.Item($destRow, 5).Value = $row[4]
}

# Copy the "A" column number/border/bold style down onto the newly written rows (A10:A17)
$ws.Cells.Item(7, 1).Copy()
$ws.Range("A10:A17").PasteSpecial(-4122)

# --- Step 3: overwrite rows 8 and 9 with the two new line entries (line7, line8) ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# --- Step 4: renumber the "A" index column sequentially for the shifted extr rows (now 10..17) ---
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Step 5: apply the refreshed from_bus / to_bus / in_service values for extr1..extr8 ---
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
